$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITR input data")
$r = 13
$ws.Range("X$r").Value = 0
$ws.Range("Y$r").Value = 0
$ws.Range("Z$r").Value = 0
$ws.Range("AA$r").Value = 0
$ws.Range("AB$r").Value = 0
$ws.Range("AE$r").Formula = "=IF(OR(ISBLANK(Q$r), ISBLANK(X$r)),`"`",Q$r+X$r)"
$ws.Range("AF$r").Formula = "=IF(ISBLANK(R$r),IF(ISBLANK(Y$r),`"`",Y$r),R$r+Y$r)"
$ws.Range("AG$r").Formula = "=IF(ISBLANK(S$r),IF(ISBLANK(Z$r),`"`",Z$r),S$r+Z$r)"
$ws.Range("AH$r").Formula = "=IF(ISBLANK(T$r),IF(ISBLANK(AA$r),`"`",AA$r),T$r+AA$r)"
$ws.Range("AI$r").Formula = "=IF(ISBLANK(U$r),IF(ISBLANK(AB$r),`"`",AB$r),U$r+AB$r)"
$ws.Range("AJ$r").Formula = "=IF(ISBLANK(V$r),IF(ISBLANK(AC$r),`"`",AC$r),V$r+AC$r)"
$ws.Range("AK$r").Formula = "=IF(ISBLANK(W$r),IF(ISBLANK(AD$r),`"`",AD$r),W$r+AD$r)"
